$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L, row 3: empty cell matching the style of row 3 (thin bottom border row)
$ws.Range("L3").Value = $null

# Column L, row 4: header year 2021, styled to match the year header row (bold, bottom border)
$ws.Range("L4").Value = 2021

# Column L, row 5: value for "Small enterprises"
$ws.Range("L5").Value = 2.3

# Column L, row 6: value for "Medium-sized enterprises"
$ws.Range("L6").Value = 1.3

# Apply formatting to match the new column with a style similar to the existing table
$ws.Range("L3").Borders.Item(3).LineStyle = 1  # bottom border continuation from row3 style
$ws.Range("L4").Font.Bold = $true
$ws.Range("L4").Borders.Item(4).LineStyle = 1
$ws.Range("L4").Borders.Item(4).Weight = -4138  # medium

$ws.Range("L6").Borders.Item(4).LineStyle = 1
$ws.Range("L6").Borders.Item(4).Weight = -4138  # medium

$ws.Cells.Item(5, 12).Select()
